$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.291003049473062
$ws.Range("C2").Value = 0.0402907006301092
$ws.Range("D2").Value = 0.6356570760909079
$ws.Range("E2").Value = 0.2523126802185374
$ws.Range("G2").Value = 0.8498863867500859
$ws.Range("H2").Value = 0.9253812374694164
$ws.Range("I2").Value = 0.7599546138712299
$ws.Range("J2").Value = 0.1262281300772727
$ws.Range("K2").Value = 0.3062434652994739
$ws.Range("M2").Value = 0.2659749248320509
$ws.Range("O2").Value = 3.576942259872723
$ws.Range("B3").Value = 0.2585444036277522
$ws.Range("C3").Value = 0.03522279083021829
$ws.Range("D3").Value = 0.630318815717871
$ws.Range("E3").Value = 0.2510763916632044
$ws.Range("G3").Value = 0.8548767174202396
$ws.Range("H3").Value = 0.9317138462302665
$ws.Range("I3").Value = 0.767041290919753
$ws.Range("J3").Value = 0.1262064540689565
$ws.Range("K3").Value = 0.270800086608773
$ws.Range("M3").Value = 0.2539989774158045
$ws.Range("O3").Value = 3.60050482414826
$ws.Range("B4").Value = 0.2386075621261341
$ws.Range("C4").Value = 0.03210142484067546
$ws.Range("D4").Value = 0.6273368094620224
$ws.Range("E4").Value = 0.2504335323714457
$ws.Range("G4").Value = 0.8584256575992057
$ws.Range("H4").Value = 0.9359624268577704
$ws.Range("I4").Value = 0.7717564610455909
$ws.Range("J4").Value = 0.1262496292161224
$ws.Range("K4").Value = 0.2490152984482421
$ws.Range("M4").Value = 0.2467448849688338
$ws.Range("O4").Value = 3.616744953434022
$ws.Range("B5").Value = 0.2304818811139171
$ws.Range("C5").Value = 0.0308270848830432
$ws.Range("D5").Value = 0.6261961423229252
$ws.Range("E5").Value = 0.2502008388517289
$ws.Range("G5").Value = 0.8599937863470188
$ws.Range("H5").Value = 0.9377844472445815
$ws.Range("I5").Value = 0.7737694535532746
$ws.Range("J5").Value = 0.1262814419288212
$ws.Range("K5").Value = 0.2401327393910293
$ws.Range("M5").Value = 0.2438139260933383
$ws.Range("O5").Value = 3.623808809085631
$ws.Range("B6").Value = 0.2291325577178043
$ws.Range("C6").Value = 0.0306153413056478
$ws.Range("D6").Value = 0.6260112414222192
$ws.Range("E6").Value = 0.2501639700123341
$ws.Range("G6").Value = 0.8602615353833514
$ws.Range("H6").Value = 0.9380924721243531
$ws.Range("I6").Value = 0.7741092374892382
$ws.Range("J6").Value = 0.1262875836110879
$ws.Range("K6").Value = 0.2386575066965833
$ws.Range("M6").Value = 0.2433287670893876
$ws.Range("O6").Value = 3.625008690338788
$ws.Range("B7").Value = 0.2384979807337118
$ws.Range("C7").Value = 0.03208424810249255
$ws.Range("D7").Value = 0.6273211240568344
$ws.Range("E7").Value = 0.2504302755768038
$ws.Range("G7").Value = 0.8584463123558876
$ws.Range("H7").Value = 0.9359866319576682
$ws.Range("I7").Value = 0.7717832383812784
$ws.Range("J7").Value = 0.1262500006633225
$ws.Range("K7").Value = 0.2488955249870202
$ws.Range("M7").Value = 0.2467052549634019
$ws.Range("O7").Value = 3.616838413574811
$ws.Range("B8").Value = 0.2798130896065629
$ws.Range("C8").Value = 0.03854532735883254
$ws.Range("D8").Value = 0.6337551667334225
$ws.Range("E8").Value = 0.2518623144609364
$ws.Range("G8").Value = 0.851506411145813
$ws.Range("H8").Value = 0.9274899820338049
$ws.Range("I8").Value = 0.7623225751478735
$ws.Range("J8").Value = 0.1262089415986694
$ws.Range("K8").Value = 0.2940275608049205
$ws.Range("M8").Value = 0.2618251481858067
$ws.Range("O8").Value = 3.584698801589553
$ws.Range("B9").Value = 0.3607562448077317
$ws.Range("C9").Value = 0.05113652641922783
$ws.Range("D9").Value = 0.6487126215383796
$ws.Range("E9").Value = 0.255590942975104
$ws.Range("G9").Value = 0.8417456022186514
$ws.Range("H9").Value = 0.9136837143313414
$ws.Range("I9").Value = 0.7466573181124332
$ws.Range("J9").Value = 0.126576078725293
$ws.Range("K9").Value = 0.3823334802848706
$ws.Range("M9").Value = 0.2922550752789519
$ws.Range("O9").Value = 3.535733677796344
$ws.Range("B10").Value = 0.4201589686708473
$ws.Range("C10").Value = 0.06033681173697403
$ws.Range("D10").Value = 0.6611225060240145
$ws.Range("E10").Value = 0.2588896698897614
$ws.Range("G10").Value = 0.8369228118106946
$ws.Range("H10").Value = 0.9052766300033852
$ws.Range("I10").Value = 0.7369074717070703
$ws.Range("J10").Value = 0.1271182156636428
$ws.Range("K10").Value = 0.447070604304173
$ws.Range("M10").Value = 0.3150802532907448
$ws.Range("O10").Value = 3.508326878810323
$ws.Range("B11").Value = 0.4471645372680655
$ws.Range("C11").Value = 0.06451086286908492
$ws.Range("D11").Value = 0.6670753961913931
$ws.Range("E11").Value = 0.260511451697063
$ws.Range("G11").Value = 0.8352393287400446
$ws.Range("H11").Value = 0.9018281410721869
$ws.Range("I11").Value = 0.7328539395505089
$ws.Range("J11").Value = 0.1274239050540515
$ws.Range("K11").Value = 0.4764865610571292
$ws.Range("M11").Value = 0.3255643513860988
$ws.Range("O11").Value = 3.497718626816123
$ws.Range("B12").Value = 0.4573879376153798
$ws.Range("C12").Value = 0.06608979936707726
$ws.Range("D12").Value = 0.6693736898942007
$ws.Range("E12").Value = 0.2611429615867777
$ws.Range("G12").Value = 0.8346752724044819
$ws.Range("H12").Value = 0.9005762742876584
$ws.Range("I12").Value = 0.73137386006238
$ws.Range("J12").Value = 0.1275481441720814
$ws.Range("K12").Value = 0.4876203568614983
$ws.Range("M12").Value = 0.3295487333502862
$ws.Range("O12").Value = 3.493968844289384
$ws.Range("B13").Value = 0.4551862889697702
$ws.Range("C13").Value = 0.06574982331792967
$ws.Range("D13").Value = 0.6688767541234029
$ws.Range("E13").Value = 0.2610061825887442
$ws.Range("G13").Value = 0.8347934849879977
$ws.Range("H13").Value = 0.9008434854029161
$ws.Range("I13").Value = 0.7316901795168391
$ws.Range("J13").Value = 0.1275210100564692
$ws.Range("K13").Value = 0.4852227439552337
$ws.Range("M13").Value = 0.3286899936485526
$ws.Range("O13").Value = 3.494764537761995
$ws.Range("B14").Value = 0.4480056864215385
$ws.Range("C14").Value = 0.06464079714392312
$ws.Range("D14").Value = 0.6672635960481728
$ws.Range("E14").Value = 0.2605630584261434
$ws.Range("G14").Value = 0.8351914513505392
$ws.Range("H14").Value = 0.9017240673060485
$ws.Range("I14").Value = 0.7327310720585949
$ws.Range("J14").Value = 0.1274339563788729
$ws.Range("K14").Value = 0.4774026565776239
$ws.Range("M14").Value = 0.3258918638949311
$ws.Range("O14").Value = 3.497404772353548
$ws.Range("B15").Value = 0.4436069494211949
$ws.Range("C15").Value = 0.06396126479648956
$ws.Range("D15").Value = 0.6662812235984461
$ws.Range("E15").Value = 0.2602938935522374
$ws.Range("G15").Value = 0.8354447828468068
$ws.Range("H15").Value = 0.9022704799743906
$ws.Range("I15").Value = 0.7333757994853087
$ws.Range("J15").Value = 0.1273817376073225
$ws.Range("K15").Value = 0.4726119076598252
$ws.Range("M15").Value = 0.3241797821316865
$ws.Range("O15").Value = 3.499056806107717
$ws.Range("B16").Value = 0.4183936995166846
$ws.Range("C16").Value = 0.06006379580558985
$ws.Range("D16").Value = 0.6607396456782624
$ws.Range("E16").Value = 0.258786116837129
$ws.Range("G16").Value = 0.8370431062200225
$ws.Range("H16").Value = 0.905509556690248
$ws.Range("I16").Value = 0.7371800619853062
$ws.Range("J16").Value = 0.127099425659523
$ws.Range("K16").Value = 0.4451474810427101
$ws.Range("M16").Value = 0.3143971030156507
$ws.Range("O16").Value = 3.509057563623088
$ws.Range("B17").Value = 0.4029214172439879
$ws.Range("C17").Value = 0.05766990166691244
$ws.Range("D17").Value = 0.6574187262513931
$ws.Range("E17").Value = 0.257892145915541
$ws.Range("G17").Value = 0.8381543897405948
$ws.Range("H17").Value = 0.9075928701628868
$ws.Range("I17").Value = 0.7396116319168691
$ws.Range("J17").Value = 0.1269413561958146
$ws.Range("K17").Value = 0.4282899906234832
$ws.Range("M17").Value = 0.3084214132515086
$ws.Range("O17").Value = 3.515668883146162
$ws.Range("B18").Value = 0.3940206011009479
$ws.Range("C18").Value = 0.05629194688300743
$ws.Range("D18").Value = 0.6555375812048396
$ws.Range("E18").Value = 0.2573893658952855
$ws.Range("G18").Value = 0.8388416120304498
$ws.Range("H18").Value = 0.9088265252761261
$ws.Range("I18").Value = 0.7410461352572142
$ws.Range("J18").Value = 0.1268559994022027
$ws.Range("K18").Value = 0.4185909108783505
$ws.Range("M18").Value = 0.3049938561139314
$ws.Range("O18").Value = 3.51964654104134
$ws.Range("B19").Value = 0.3910066857909555
$ws.Range("C19").Value = 0.05582521698097764
$ws.Range("D19").Value = 0.6549056367398691
$ws.Range("E19").Value = 0.2572210940977939
$ws.Range("G19").Value = 0.8390825435219114
$ws.Range("H19").Value = 0.909250299450818
$ws.Range("I19").Value = 0.7415380037906552
$ws.Range("J19").Value = 0.1268280544834113
$ws.Range("K19").Value = 0.4153064555733863
$ws.Range("M19").Value = 0.3038349826105531
$ws.Range("O19").Value = 3.521023365106942
$ws.Range("B20").Value = 0.4045686346113371
$ws.Range("C20").Value = 0.05792484501829165
$ws.Range("D20").Value = 0.6577692477406458
$ws.Range("E20").Value = 0.2579861302693445
$ws.Range("G20").Value = 0.8380311192927365
$ws.Range("H20").Value = 0.9073674357506576
$ws.Range("I20").Value = 0.7393490685661739
$ws.Range("J20").Value = 0.126957607553237
$ws.Range("K20").Value = 0.4300848241643962
$ws.Range("M20").Value = 0.3090565538088512
$ws.Range("O20").Value = 3.514946985037881
$ws.Range("B21").Value = 0.4501148900409362
$ws.Range("C21").Value = 0.06496659136946903
$ws.Range("D21").Value = 0.6677362251280101
$ws.Range("E21").Value = 0.260692743532303
$ws.Range("G21").Value = 0.835072565442843
$ws.Range("H21").Value = 0.901463953940862
$ws.Range("I21").Value = 0.7324238464652311
$ws.Range("J21").Value = 0.1274592960794223
$ws.Range("K21").Value = 0.4796997568379879
$ws.Range("M21").Value = 0.3267133565456888
$ws.Range("O21").Value = 3.496622016326512
$ws.Range("B22").Value = 0.4798641190267574
$ws.Range("C22").Value = 0.06955891738581954
$ws.Range("D22").Value = 0.6745069793460061
$ws.Range("E22").Value = 0.2625629316878815
$ws.Range("G22").Value = 0.8335670762063643
$ws.Range("H22").Value = 0.8979204205864022
$ws.Range("I22").Value = 0.7282178493289209
$ws.Range("J22").Value = 0.1278366073856887
$ws.Range("K22").Value = 0.5120942963067989
$ws.Range("M22").Value = 0.3383362420089711
$ws.Range("O22").Value = 3.486203795472534
$ws.Range("B23").Value = 0.4639882209776829
$ws.Range("C23").Value = 0.06710883362454467
$ws.Range("D23").Value = 0.6708698616046433
$ws.Range("E23").Value = 0.2615555278072321
$ws.Range("G23").Value = 0.8343313991682351
$ws.Range("H23").Value = 0.8997828919765283
$ws.Range("I23").Value = 0.7304333812222339
$ws.Range("J23").Value = 0.1276307107343726
$ws.Range("K23").Value = 0.494807818826672
$ws.Range("M23").Value = 0.3321253525550105
$ws.Range("O23").Value = 3.491621628179075
$ws.Range("B24").Value = 0.403823945024385
$ws.Range("C24").Value = 0.05780959033575073
$ws.Range("D24").Value = 0.6576106894494274
$ws.Range("E24").Value = 0.2579436051284318
$ws.Range("G24").Value = 0.8380866993184668
$ws.Range("H24").Value = 0.9074692427220867
$ws.Range("I24").Value = 0.7394676595710337
$ws.Range("J24").Value = 0.1269502431214562
$ws.Range("K24").Value = 0.4292734031874375
$ws.Range("M24").Value = 0.3087693819948996
$ws.Range("O24").Value = 3.515272804689204
$ws.Range("B25").Value = 0.3388691593967508
$ws.Range("C25").Value = 0.04773896876962169
$ws.Range("D25").Value = 0.6444163140087369
$ws.Range("E25").Value = 0.2544839003010217
$ws.Range("G25").Value = 0.8439739170653695
$ws.Range("H25").Value = 0.9171134260457166
$ws.Range("I25").Value = 0.7505861685231885
$ws.Range("J25").Value = 0.1264288857082008
$ws.Range("K25").Value = 0.3584677100649856
$ws.Range("M25").Value = 0.2839402202306474
$ws.Range("O25").Value = 3.547475155863424
